$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.178117048346056
$ws.Range("C2").Value = 0.5903307888040712
$ws.Range("J2").Value = 0.01272264631043257
$ws.Range("P2").Value = 0.1094147582697201
$ws.Range("S2").Value = 0.1094147582697201
# Row 3
$ws.Range("B3").Value = 0.00819672131147541
$ws.Range("C3").Value = 0.03688524590163934
$ws.Range("J3").Value = 0.02459016393442623
$ws.Range("P3").Value = 0.75
$ws.Range("S3").Value = 0.180327868852459
# Row 4
$ws.Range("J4").Value = 0.05
$ws.Range("P4").Value = 0.8
$ws.Range("S4").Value = 0.15
# Row 6
$ws.Range("B6").Value = 0.0502283105022831
$ws.Range("D6").Value = 0.0182648401826484
$ws.Range("F6").Value = 0.0502283105022831
$ws.Range("J6").Value = 0.3150684931506849
$ws.Range("O6").Value = 0.0045662100456621
$ws.Range("Q6").Value = 0.1461187214611872
$ws.Range("R6").Value = 0.0502283105022831
$ws.Range("S6").Value = 0.365296803652968
# Row 7
$ws.Range("B7").Value = 0.1036036036036036
$ws.Range("D7").Value = 0.02702702702702703
$ws.Range("E7").Value = 0.004504504504504504
$ws.Range("F7").Value = 0.04504504504504504
$ws.Range("J7").Value = 0.1441441441441441
$ws.Range("O7").Value = 0.01351351351351351
$ws.Range("Q7").Value = 0.2162162162162162
$ws.Range("R7").Value = 0.0945945945945946
$ws.Range("S7").Value = 0.3513513513513514
# Row 8
$ws.Range("B8").Value = 0.125968992248062
$ws.Range("D8").Value = 0.01937984496124031
$ws.Range("F8").Value = 0.07364341085271318
$ws.Range("J8").Value = 0.1007751937984496
$ws.Range("O8").Value = 0.02713178294573643
$ws.Range("Q8").Value = 0.186046511627907
$ws.Range("R8").Value = 0.05813953488372093
$ws.Range("S8").Value = 0.4089147286821705
# Row 9
$ws.Range("B9").Value = 0.1693989071038251
$ws.Range("D9").Value = 0.01092896174863388
$ws.Range("E9").Value = 0.00546448087431694
$ws.Range("F9").Value = 0.03825136612021858
$ws.Range("J9").Value = 0.07103825136612021
$ws.Range("O9").Value = 0.02185792349726776
$ws.Range("Q9").Value = 0.1584699453551913
$ws.Range("R9").Value = 0.0546448087431694
$ws.Range("S9").Value = 0.4699453551912569
# Row 10
$ws.Range("B10").Value = 0.1476190476190476
$ws.Range("D10").Value = 0.01825396825396826
$ws.Range("E10").Value = 0.001587301587301587
$ws.Range("F10").Value = 0.0626984126984127
$ws.Range("J10").Value = 0.1007936507936508
$ws.Range("O10").Value = 0.01349206349206349
$ws.Range("Q10").Value = 0.2119047619047619
$ws.Range("R10").Value = 0.08888888888888889
$ws.Range("S10").Value = 0.3547619047619048
# Row 11
$ws.Range("G11").Value = 0.1798365122615804
$ws.Range("J11").Value = 0.09809264305177112
$ws.Range("K11").Value = 0.2343324250681199
$ws.Range("L11").Value = 0.4713896457765668
$ws.Range("S11").Value = 0.01634877384196185
# Row 12
$ws.Range("G12").Value = 0.6830601092896175
$ws.Range("J12").Value = 0.2513661202185792
$ws.Range("K12").Value = 0.02185792349726776
$ws.Range("L12").Value = 0.01639344262295082
$ws.Range("S12").Value = 0.0273224043715847
# Row 13
$ws.Range("F13").Value = 0.01470588235294118
$ws.Range("G13").Value = 0.5588235294117647
$ws.Range("J13").Value = 0.3529411764705883
$ws.Range("S13").Value = 0.07352941176470588
# Row 15
$ws.Range("F15").Value = 0.01762114537444934
$ws.Range("H15").Value = 0.1718061674008811
$ws.Range("I15").Value = 0.04405286343612335
$ws.Range("J15").Value = 0.3171806167400881
$ws.Range("K15").Value = 0.09691629955947137
$ws.Range("M15").Value = 0.00881057268722467
$ws.Range("O15").Value = 0.09691629955947137
$ws.Range("S15").Value = 0.2466960352422908
# Row 16
$ws.Range("F16").Value = 0.024
$ws.Range("H16").Value = 0.172
$ws.Range("I16").Value = 0.112
$ws.Range("J16").Value = 0.364
$ws.Range("K16").Value = 0.116
$ws.Range("M16").Value = 0.02
$ws.Range("O16").Value = 0.064
$ws.Range("S16").Value = 0.128
# Row 17
$ws.Range("F17").Value = 0.02748414376321353
$ws.Range("H17").Value = 0.1733615221987315
$ws.Range("I17").Value = 0.09513742071881606
$ws.Range("J17").Value = 0.4038054968287527
$ws.Range("K17").Value = 0.09936575052854123
$ws.Range("M17").Value = 0.02325581395348837
$ws.Range("N17").Value = 0.002114164904862579
$ws.Range("O17").Value = 0.06553911205073996
$ws.Range("S17").Value = 0.1099365750528541
# Row 18
$ws.Range("F18").Value = 0.01685393258426966
$ws.Range("H18").Value = 0.2247191011235955
$ws.Range("I18").Value = 0.07303370786516854
$ws.Range("J18").Value = 0.3764044943820224
$ws.Range("K18").Value = 0.1179775280898876
$ws.Range("M18").Value = 0.02808988764044944
$ws.Range("O18").Value = 0.06179775280898876
$ws.Range("S18").Value = 0.101123595505618
# Row 19
$ws.Range("F19").Value = 0.01470588235294118
$ws.Range("H19").Value = 0.2438080495356037
$ws.Range("I19").Value = 0.0696594427244582
$ws.Range("J19").Value = 0.3459752321981424
$ws.Range("K19").Value = 0.1207430340557276
$ws.Range("M19").Value = 0.03869969040247678
$ws.Range("N19").Value = 0.0007739938080495357
$ws.Range("O19").Value = 0.06191950464396285
$ws.Range("S19").Value = 0.1037151702786378
